$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-12-16 08:25:52"

$wsZhCn.Range("H4").Value = "2016-12-16 08:25:39"
$wsZhCn.Range("L4").Value = "2016-12-16 08:26:35"

$wsDeDe.Range("H4").Value = "2016-12-16 08:25:52"
$wsDeDe.Range("L4").Value = "2016-12-16 08:26:53"
